$wb = $excel.ActiveWorkbook

# Rename worksheets: rfStability -> rsStability, rfRank -> rsRank
$wsStability = $wb.Worksheets.Item("rfStability")
$wsStability.Name = "rsStability"

$wsRank = $wb.Worksheets.Item("rfRank")
$wsRank.Name = "rsRank"

# Update header row text on both sheets (headers are shared across sheets)
foreach ($ws in @($wsStability, $wsRank)) {
    $ws.Range("A1").Value = "Target"
    $ws.Range("B1").Value = "delta-Ct"
    $ws.Range("C1").Value = "BestKeeper"
    $ws.Range("D1").Value = "Normfinder"
    $ws.Range("E1").Value = "geNorm"
    $ws.Range("F1").Value = "Comprehensive Rank"
}
